$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add "username" header in A1, matching the style of the
#     existing header cells (bold / bordered / centered = same style as B1) ---
$ws.Range("A1").Value = "username"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 2: existing user "sumona" keeps her name; password + balance change ---
$ws.Range("B2").Value = "'123"        # force text (password column is unformatted)
$ws.Range("C2").Value = 988

# --- Row 3: new user "sarmin" ---
$ws.Range("A3").Value = "sarmin"
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)   # xlPasteFormats (same style as other usernames)
$ws.Range("B3").Value = "nan"         # not numeric-looking, no text-forcing needed
$ws.Range("C3").Value = 1012

# --- Row 4: new user "ria" ---
$ws.Range("A4").Value = "ria"
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B4").Value = "'89"
$ws.Range("C4").Value = 1000

# Strip the transient quote-prefix ("forced text") formatting the apostrophe
# trick applied above so the password cells end up unformatted again, same
# as B2 originally was.
$ws.Range("B2").ClearFormats()
$ws.Range("B4").ClearFormats()
